$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.372.95"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -1.12%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.837.49"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +2.53%  "

$ws.Range("E4").Value = "  -0.06%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.76"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.44%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.07"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -3.10%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.833.39"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +2.51%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -2.65%  "

$ws.Range("E10").Value = "  -2.17%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -1.30%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -0.89%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.68"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -3.53%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -2.21%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.474.61"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +2.45%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.837.53"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +2.46%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.572.20"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("E18").Value = "  +1.82%  "

$ws.Range("E19").Value = "  -0.57%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.07"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -1.57%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.25"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.87%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.23"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -1.70%  "

$ws.Range("E24").Value = "  +5.87%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.79"
$ws.Range("D25").Style = $style

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -3.35%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.06"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -1.52%  "

$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("E29").Value = "  -1.88%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.94"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -1.42%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.987.36"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +2.56%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.81"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -4.49%  "

$ws.Range("E33").Value = "  -4.46%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.95"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +1.45%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.781.44"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +2.91%  "

$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("E37").Value = "  +0.86%  "

$ws.Range("E38").Value = "  +0.04%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.84"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -1.81%  "

$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("E41").Value = "  -2.82%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.94"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -3.35%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "428.03"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +1.03%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.49"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("E45").Value = "  -1.22%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.34"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -1.35%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "142.84"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +0.57%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.827.73"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +1.65%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.10"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +14.32%  "

$ws.Range("E51").Value = "  +0.30%  "
